$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.027.82'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '2.244.45'
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '98.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +17.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '270.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.27%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.641'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0944'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +17.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("E14").Value = '  +7.27%  '
$ws.Range("D15").Value = '2.572.52'
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("E16").Value = '  +5.24%  '
$ws.Range("D17").Value = '2.236.99'
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").Value = '44.011.70'
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("E20").Value = '  +5.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.49%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.43%  '
$ws.Range("E28").Value = '  +2.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0928'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.26%  '
$ws.Range("E34").Value = '  +5.11%  '
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.70'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +31.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.248'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +24.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("E42").Value = '  +4.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.102'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("E48").Value = '  +4.52%  '
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.460.74'
$ws.Range("E51").Value = '  +2.19%  '
